$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes old row 3 down to row 4),
# so we end up with 4 data rows (2-5) splitting the old "M2" target
# cluster into separate "M1" and "M2" rows per sending cluster.
$ws.Rows.Item(3).Insert()

$data = New-Object 'object[,]' 4,20
# Row 2: ECs / Pomc / Oprm1 / M1
$data[0,0] = "ECs"
$data[0,1] = "Pomc"
$data[0,2] = "Oprm1"
$data[0,3] = "M1"
$data[0,4] = 2
$data[0,5] = 1
$data[0,6] = 1.4951615
$data[0,7] = 2.990323
$data[0,8] = 0.6020739711267923
$data[0,9] = 0.5021622551131893
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.05215766666666666
$data[0,13] = 0.156473
$data[0,14] = 0.1010355835763341
$data[0,15] = 0.1010355835763341
$data[0,16] = 0.07798413512983333
$data[0,17] = 0.467904810779
$data[0,18] = 0.06083089502891637
$data[0,19] = 0.05073625649536904
# Row 3: ECs / Pomc / Oprm1 / M2
$data[1,0] = "ECs"
$data[1,1] = "Pomc"
$data[1,2] = "Oprm1"
$data[1,3] = "M2"
$data[1,4] = 2
$data[1,5] = 1
$data[1,6] = 1.4951615
$data[1,7] = 2.990323
$data[1,8] = 0.6020739711267923
$data[1,9] = 0.5021622551131893
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.464073
$data[1,13] = 1.392219
$data[1,14] = 0.8989644164236659
$data[1,15] = 0.8989644164236659
$data[1,16] = 0.6938640827895001
$data[1,17] = 4.163184496737
$data[1,18] = 0.5412430760978759
$data[1,19] = 0.4514259986178203
# Row 4: Neutro / Pomc / Oprm1 / M1
$data[2,0] = "Neutro"
$data[2,1] = "Pomc"
$data[2,2] = "Oprm1"
$data[2,3] = "M1"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 0.9881903333333333
$data[2,7] = 2.964571
$data[2,8] = 0.3979260288732077
$data[2,9] = 0.4978377448868108
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.05215766666666666
$data[2,13] = 0.156473
$data[2,14] = 0.1010355835763341
$data[2,15] = 0.1010355835763341
$data[2,16] = 0.05154170200922222
$data[2,17] = 0.463875318083
$data[2,18] = 0.0402046885474177
$data[2,19] = 0.05029932708096505
# Row 5: Neutro / Pomc / Oprm1 / M2
$data[3,0] = "Neutro"
$data[3,1] = "Pomc"
$data[3,2] = "Oprm1"
$data[3,3] = "M2"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 0.9881903333333333
$data[3,7] = 2.964571
$data[3,8] = 0.3979260288732077
$data[3,9] = 0.4978377448868108
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.464073
$data[3,13] = 1.392219
$data[3,14] = 0.8989644164236659
$data[3,15] = 0.8989644164236659
$data[3,16] = 0.458592452561
$data[3,17] = 4.127332073049
$data[3,18] = 0.35772134032579
$data[3,19] = 0.4475384178058457

$ws.Range("A2:T5").Value = $data

Write-Output ("Dimension should now be A1:T5; UsedRange: " + $ws.UsedRange.Address())
